$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment") entirely, shifting columns C:J left to B:I.
$ws.Columns.Item(2).Delete()

# Append ".jamais.jamais" to each header title in row 1 (columns B through I).
for ($col = 2; $col -le 9; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = "$($cell.Value2).jamais.jamais"
}
